# feat: added personnel reassignment page
# Adds a new "instrumental case" column (E) to the "Ролі" (Roles) sheet
# and makes that sheet the active tab/selection, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column header + four data rows (instrumental case of each role name)
$ws.Range("E1").Value = "Назва ролі в орудному відмінку"
$ws.Range("E2").Value = "старшим науковим співробітником"
$ws.Range("E3").Value = "дідоводом"
$ws.Range("E4").Value = "водієм - електриком"
$ws.Range("E5").Value = "начальником служби"

# Switch the active tab to the Roles sheet and leave the selection on E5,
# mirroring the workbook/sheet-view changes in the target file.
$ws.Activate()
$ws.Range("E5").Select()
